$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Sheet1 - scene/quest data
$ws2 = $wb.Worksheets.Item(2)   # MapSet - npc/map data

# ---------------------------------------------------------------
# Sheet1: insert a new row 27 for the "hiddeway" (hidden passage)
# scene entry, pushing the former rows 27-29 down to 28-30.
# ---------------------------------------------------------------
$ws1.Rows.Item(27).Insert()

$ws1.Range("A27").Value = 42010018
$ws1.Range("F27").Value = "hiddeway"
$ws1.Range("B27").Value = "隐藏通道"
$ws1.Range("C27").Value = 1
$ws1.Range("D27").Value = 0
$ws1.Range("G27").Value = "hiddeway"
$ws1.Range("H27").Value = "hiddeway"

# Grow table1 (表3) so it keeps covering the data, now through row 30.
$lo1 = $ws1.ListObjects.Item(1)
$lo1.Resize($ws1.Range("A3:AC30"))

# ---------------------------------------------------------------
# Sheet2 (MapSet): point row 11 (Aolai) at its replacement via the
# new ReplaceId column, and append a new row 18 describing the
# replacement npc "npcaolai2".
# ---------------------------------------------------------------
$ws2.Range("U11").Value = 42030012

$ws2.Range("A18").Value = 42030012
$ws2.Range("F18").Value = "npcaolai2"
$ws2.Range("B18").Value = "奥莱伊李"
$ws2.Range("C18").Value = 2
$ws2.Range("D18").Value = 0
$ws2.Range("G18").Value = "npcaolai"
$ws2.Range("H18").Value = "npcaolai2"
$ws2.Range("I18").Value = "true"
$ws2.Range("J18").Value = 43020105
$ws2.Range("K18").Value = "oneline"
$ws2.Range("S18").Value = 43020105

# Grow table2 (表3_5) so it keeps covering the data, now through row 18.
$lo2 = $ws2.ListObjects.Item(1)
$lo2.Resize($ws2.Range("A3:AC18"))

# ---------------------------------------------------------------
# Page setup for Sheet1 (A4 portrait), matching the finished sheet.
# ---------------------------------------------------------------
$ws1.PageSetup.PaperSize = 9
$ws1.PageSetup.Orientation = 1

# ---------------------------------------------------------------
# Restore the selections left on each sheet, then make Sheet1 the
# active (selected) tab, as it ends up in the final workbook.
# ---------------------------------------------------------------
$ws2.Range("U11").Select()

$ws1.Activate()
$ws1.Range("U4").Select()
